$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove existing rows 1-3 entirely so stale row-height / thick-border
# flags (ht=, thickBot=) from the old file do not linger in the output ---
$ws.Rows("1:3").EntireRow.Delete()

# --- Row 1 (header) ---
$ws.Range("A1").Value = "lang_code"
$ws.Range("B1").Value = "code"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "descr"
$ws.Range("E1").Value = "is_active"

$hdr = $ws.Range("A1:D1")
$hdr.Font.Name = "Cambria"
$hdr.Font.Bold = $true
$hdr.Font.Color = 0
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.WrapText = $false
$hdr.Borders.LineStyle = 1
$hdr.Borders.Weight = 2

$hdrE = $ws.Range("E1")
$hdrE.Font.Name = "Cambria"
$hdrE.Font.Bold = $true
$hdrE.Font.Color = 0
$hdrE.HorizontalAlignment = -4108
$hdrE.VerticalAlignment = -4160
$hdrE.WrapText = $false
$hdrE.Borders.LineStyle = 1
$hdrE.Borders.Weight = 2
$hdrE.NumberFormat = "@"

# --- Row 2 : MNA / Manual Adjudication reason category ---
$ws.Range("A2").Value = "eng"
$ws.Range("B2").Value = "MNA"
$ws.Range("C2").Value = "Manual Adjudication"
$ws.Range("D2").Value = "Rejection during Manual Adjudication"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "'TRUE"

# --- Row 3 : CLR / Client Rejection reason category ---
$ws.Range("A3").Value = "eng"
$ws.Range("B3").Value = "CLR"
$ws.Range("C3").Value = "Client Rejection"
$ws.Range("D3").Value = "Rejection in Registration Client"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "'TRUE"

# --- New blank, formatted rows 6 & 7 ---
$blank = $ws.Range("C6:D7")
$blank.HorizontalAlignment = -4131
$blank.WrapText = $true

# --- Column widths (closest achievable to target 31.26953125 / 46.1796875
# given this engine's internal character-width rounding granularity) ---
$ws.Columns("C").ColumnWidth = 30.417
$ws.Columns("D").ColumnWidth = 45.25

# --- Selection / view ---
$ws.Range("E14").Select()

# --- Page margins (inches) to match target (converted from the cm-based
# values baked into the workbook) ---
$ws.PageSetup.LeftMargin = $excel.CentimetersToPoints(1.905)
$ws.PageSetup.RightMargin = $excel.CentimetersToPoints(1.905)
$ws.PageSetup.TopMargin = $excel.CentimetersToPoints(2.54)
$ws.PageSetup.BottomMargin = $excel.CentimetersToPoints(2.54)
$ws.PageSetup.HeaderMargin = $excel.CentimetersToPoints(1.3)
$ws.PageSetup.FooterMargin = $excel.CentimetersToPoints(1.3)
